$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up the header row capitalization / spacing so the columns read
# "Students", "test 1", "Test 2", "Test 3", "Midterm", "Test 4", "Test 5",
# "Test 6", "Final".
$ws.Range("A1").Value = "Students"
$ws.Range("B1").Value = "test 1"
$ws.Range("C1").Value = "Test 2"
$ws.Range("D1").Value = "Test 3"
$ws.Range("E1").Value = "Midterm"
$ws.Range("F1").Value = "Test 4"
$ws.Range("G1").Value = "Test 5"
$ws.Range("H1").Value = "Test 6"
$ws.Range("I1").Value = "Final"

# rathi.kashi_ug20 didn't take the latest test ("Final") — remove the
# stray score that had been entered in I3.
$ws.Range("I3").ClearContents()

# Leave the selection on the cell that was just cleared.
$ws.Range("I3").Select()
